$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '43.325.01'
$ws.Range("E2").Value = '  -0.67%  '
$ws.Range("D3").Value = '2.277.87'
$ws.Range("E3").Value = '  -0.14%  '
$ws.Range("E4").Value = '  -0.44%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '111.79'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '264.14'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.97%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.641'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +3.10%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.608'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '46.57'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0936'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.25'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +5.27%  '
$ws.Range("E13").Value = '  +1.70%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.33'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.47%  '
$ws.Range("D15").Value = '2.620.14'
$ws.Range("E15").Value = '  -0.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.860'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.56%  '
$ws.Range("D17").Value = '2.278.81'
$ws.Range("E17").Value = '  +0.00%  '
$ws.Range("D18").Value = '43.172.02'
$ws.Range("E18").Value = '  -0.74%  '
$ws.Range("E19").Value = '  -0.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.74'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.74%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.20'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.44'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.40'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.56%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.86'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +3.88%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.36'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.44%  '
$ws.Range("E26").Value = '  +1.96%  '
$ws.Range("E27").Value = '  -1.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '41.27'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.44%  '
$ws.Range("E29").Value = '  -1.41%  '
$ws.Range("E30").Value = '  -0.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.41'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.49'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.54%  '
$ws.Range("E33").Value = '  -2.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.66'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.89%  '
$ws.Range("E35").Value = '  +3.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0381'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +6.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.69'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.86'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +3.70%  '
$ws.Range("E39").Value = '  -2.98%  '
$ws.Range("E40").Value = '  +8.49%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.27'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +5.99%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '75.43'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +6.52%  '
$ws.Range("E43").Value = '  -1.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.10'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("E46").Value = '  -0.91%  '
$ws.Range("E47").Value = '  +4.63%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.54'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0991'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '100.37'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.434'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.89%  '
